# Updates cryptocurrency price (D) and volume-change (E) values
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.947.58'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.638.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  +0.86%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0638'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0795'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.864.85'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.632.10'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0758'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.57'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.964.60'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '  +0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E23').Value = '  -1.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('E27').Value = '  +2.56%  '
$ws.Range('E28').Value = '  -0.78%  '
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('E30').Value = '  -0.38%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.30'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.40%  '
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.94%  '
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.904'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.139.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.546'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  +0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.30'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.774.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('E45').Value = '  +6.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0530'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('E49').Value = '  -0.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.63'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0963'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.08%  '
